$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily Updates")

# Semantic content change: header I1 renamed from CREATED_DATE to DUE_DATE
$ws.Range("I1").Value = "DUE_DATE"

# Update selection / scroll position to match the saved view state
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("D2").Select() | Out-Null
